# Update header names in the "QuantitySI" sheet to be SciCat compliant:
#   si_value -> valueSI
#   si_unit  -> unitSI

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitySI")

$ws.Range("A1").Value = "valueSI"
$ws.Range("B1").Value = "unitSI"
